$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.108.80"
$ws.Range("E2").Value = "  +6.63%  "

# Row 3
$ws.Range("D3").Value = "3.017.68"
$ws.Range("E3").Value = "  +3.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +13.05%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "3.014.15"
$ws.Range("E8").Value = "  +3.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.75%  "

# Row 10
$ws.Range("E10").Value = "  -2.74%  "

# Row 11
$ws.Range("E11").Value = "  +5.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.73%  "

# Row 13
$ws.Range("E13").Value = "  +8.22%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.83"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.10%  "

# Row 15
$ws.Range("E15").Value = "  -0.56%  "

# Row 16
$ws.Range("D16").Value = "66.048.08"
$ws.Range("E16").Value = "  +6.60%  "

# Row 17
$ws.Range("D17").Value = "3.518.38"
$ws.Range("E17").Value = "  +3.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +7.04%  "

# Row 19
$ws.Range("D19").Value = "3.015.47"
$ws.Range("E19").Value = "  +3.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.691"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +8.07%  "

# Row 24
$ws.Range("E24").Value = "  +4.56%  "

# Row 25
$ws.Range("E25").Value = "  +13.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.60%  "

# Row 27
$ws.Range("E27").Value = "  +3.60%  "

# Row 28
$ws.Range("E28").Value = "  -0.10%  "

# Row 29
$ws.Range("E29").Value = "  +16.33%  "

# Row 30
$ws.Range("E30").Value = "  +17.56%  "

# Row 31
$ws.Range("E31").Value = "  -6.63%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.37"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.74%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.84%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.995"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.37%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.85"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.43%  "

# Row 38
$ws.Range("E38").Value = "  +15.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.90"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.19%  "

# Row 41
$ws.Range("E41").Value = "  +16.60%  "

# Row 42
$ws.Range("E42").Value = "  +7.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.36"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.48"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.40"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +13.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0361"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.40%  "

# Row 47
$ws.Range("D47").Value = "2.807.12"
$ws.Range("E47").Value = "  +3.40%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.55"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.46%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.99"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.95%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +11.40%  "
